$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 6193.5
$ws.Range("I45").Value = 3499
$ws.Range("J45").Value = 8888
$ws.Range("K45").Value = 10497
$ws.Range("L45").Value = 26664
$ws.Range("M45").Value = -10305
$ws.Range("N45").Value = -27048

$ws.Range("H46").Value = 3832.3333
$ws.Range("I46").Value = 3748.5
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 11245.5
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = -11126.5
$ws.Range("N46").Value = -12238

$ws.Range("H49").Value = 616.3333
$ws.Range("I49").Value = 616.3333
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 1848.9999
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H58").Value = 214.4
$ws.Range("I58").Value = 237.22223
$ws.Range("J58").Value = 9
$ws.Range("K58").Value = 711.66669
$ws.Range("L58").Value = 27
$ws.Range("M58").Value = -561.66669
$ws.Range("N58").Value = -327

$ws.Range("H59").Value = 3499
$ws.Range("I59").Value = 3499
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 10497
$ws.Range("L59").ClearContents()
$ws.Range("M59").Value = -9940

$ws.Range("H60").Value = 3832.3333
$ws.Range("I60").Value = 3748.5
$ws.Range("J60").Value = 4000
$ws.Range("K60").Value = 11245.5
$ws.Range("L60").Value = 12000
$ws.Range("M60").Value = -10761.5
$ws.Range("N60").Value = -12968

$ws.Range("H92").Value = 84029.164
$ws.Range("I92").Value = 91531.82000000001
$ws.Range("J92").Value = 1500
$ws.Range("K92").Value = 91531.82000000001
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = -90283.82000000001
$ws.Range("N92").Value = -3996

$ws.Range("H113").Value = 3999.75
$ws.Range("I113").Value = 3999.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3999.75
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H116").Value = 5500
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 7000
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = -558
$ws.Range("N116").Value = -13884

$ws.Range("H136").Value = 174666.67
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 174666.67
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 174666.67
$ws.Range("N136").Value = -184866.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8138.6875
$ws.Range("I32").Value = 8078.6772
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 8078.6772
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -7791.6772
$ws.Range("N32").Value = -10573

$ws.Range("H61").Value = 4892.3
$ws.Range("I61").Value = 4115.5
$ws.Range("J61").Value = 7999.5
$ws.Range("K61").Value = 4115.5
$ws.Range("L61").Value = 7999.5
$ws.Range("M61").Value = -3903.5
$ws.Range("N61").Value = -8423.5

$ws.Range("H74").Value = 2519.6428
$ws.Range("I74").Value = 2525.4546
$ws.Range("J74").Value = 2498.3333
$ws.Range("K74").Value = 2525.4546
$ws.Range("L74").Value = 2498.3333
$ws.Range("M74").Value = -1651.4546
$ws.Range("N74").Value = -4246.3333

$ws.Range("H77").Value = 2519.6428
$ws.Range("I77").Value = 2525.4546
$ws.Range("J77").Value = 2498.3333
$ws.Range("K77").Value = 12627.273
$ws.Range("L77").Value = 12491.6665
$ws.Range("M77").Value = -8259.273000000001
$ws.Range("N77").Value = -21227.6665

$ws.Range("H88").Value = 2788.5715
$ws.Range("I88").Value = 2390.5
$ws.Range("J88").Value = 2947.8
$ws.Range("K88").Value = 2390.5
$ws.Range("L88").Value = 2947.8
$ws.Range("M88").Value = -1984.5
$ws.Range("N88").Value = -3759.8

$ws.Range("H91").Value = 2788.5715
$ws.Range("I91").Value = 2390.5
$ws.Range("J91").Value = 2947.8
$ws.Range("K91").Value = 2390.5
$ws.Range("L91").Value = 2947.8
$ws.Range("M91").Value = -986.5
$ws.Range("N91").Value = -5755.8

$ws.Range("H102").Value = 4005
$ws.Range("I102").Value = 4010
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 4010
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -2388
$ws.Range("N102").Value = -7244

$ws.Range("H110").Value = 9843.777
$ws.Range("I110").Value = 10199.25
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 10199.25
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = -8154.25
$ws.Range("N110").Value = -11090

$ws.Range("H122").Value = 2449.92
$ws.Range("I122").Value = 1972.35
$ws.Range("J122").Value = 4360.2
$ws.Range("K122").Value = 5917.049999999999
$ws.Range("L122").Value = 13080.6
$ws.Range("M122").Value = -3467.049999999999
$ws.Range("N122").Value = -17980.6

$ws.Range("H132").Value = 3430.1428
$ws.Range("I132").Value = 3430.1428
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10290.4284
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7760.428400000001

$ws.Range("H136").Value = 4892.3
$ws.Range("I136").Value = 4115.5
$ws.Range("J136").Value = 7999.5
$ws.Range("K136").Value = 12346.5
$ws.Range("L136").Value = 23998.5
$ws.Range("M136").Value = -9796.5
$ws.Range("N136").Value = -29098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6864.8887
$ws.Range("I86").Value = 2840.7144
$ws.Range("J86").Value = 20949.5
$ws.Range("K86").Value = 2840.7144
$ws.Range("L86").Value = 20949.5
$ws.Range("M86").Value = -1717.7144
$ws.Range("N86").Value = -23195.5

$ws.Range("H89").Value = 6864.8887
$ws.Range("I89").Value = 2840.7144
$ws.Range("J89").Value = 20949.5
$ws.Range("K89").Value = 14203.572
$ws.Range("L89").Value = 104747.5
$ws.Range("M89").Value = -8587.572
$ws.Range("N89").Value = -115979.5

$ws.Range("H134").Value = 7929.4
$ws.Range("I134").Value = 9268.571
$ws.Range("J134").Value = 4804.6665
$ws.Range("K134").Value = 27805.713
$ws.Range("L134").Value = 14413.9995
$ws.Range("M134").Value = -25270.713
$ws.Range("N134").Value = -19483.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 7816.6665
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 7816.6665
$ws.Range("K14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("N14").Value = -8156.6665

$ws.Range("H99").Value = 2988.2
$ws.Range("I99").Value = 2988.2
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2988.2
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1490.2

$ws.Range("H107").Value = 708.4167
$ws.Range("I107").Value = 498.9091
$ws.Range("J107").Value = 3013
$ws.Range("K107").Value = 498.9091
$ws.Range("L107").Value = 3013
$ws.Range("M107").Value = 1421.0909
$ws.Range("N107").Value = -6853

$ws.Range("H122").Value = 5610.8887
$ws.Range("I122").Value = 6212
$ws.Range("J122").Value = 3507
$ws.Range("K122").Value = 18636
$ws.Range("L122").Value = 10521
$ws.Range("M122").Value = -16186
$ws.Range("N122").Value = -15421

$ws.Range("H126").Value = 2988.2
$ws.Range("I126").Value = 2988.2
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8964.599999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6494.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 413.72726
$ws.Range("I14").Value = 413.72726
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1241.18178
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1068.18178

$ws.Range("H21").Value = 400
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("N21").Value = -1546

$ws.Range("H39").Value = 5428.5713
$ws.Range("I39").Value = 3003
$ws.Range("J39").Value = 5615.154
$ws.Range("K39").Value = 9009
$ws.Range("L39").Value = 16845.462
$ws.Range("M39").Value = -8715
$ws.Range("N39").Value = -17433.462

$ws.Range("H52").Value = 6200
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 6200
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 18600
$ws.Range("N52").Value = -19132

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1724.25
$ws.Range("I40").Value = 1724.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1724.25
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1588.25

$ws.Range("H93").Value = 1591.5
$ws.Range("I93").Value = 1824.75
$ws.Range("J93").Value = 1125
$ws.Range("K93").Value = 1824.75
$ws.Range("L93").Value = 1125
$ws.Range("M93").Value = -576.75
$ws.Range("N93").Value = -3621

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").ClearContents()

$ws.Range("H132").Value = 2848.5
$ws.Range("I132").Value = 2848.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8545.5
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 16690609
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 16690609
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 16690609
$ws.Range("N46").Value = -16691071

$ws.Range("H104").Value = 24650
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 24650
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 24650
$ws.Range("N104").Value = -31638

$ws.Range("H107").Value = 496.33334
$ws.Range("I107").Value = 494.5
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 1483.5
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 436.5
$ws.Range("N107").Value = -5340

$ws.Range("H122").Value = 8399.799999999999
$ws.Range("I122").Value = 8399.799999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 25199.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -22749.4

$ws.Range("H126").Value = 4619.2856
$ws.Range("I126").Value = 4222.5
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 12667.5
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -10197.5
$ws.Range("N126").Value = -25940

$ws.Range("H132").Value = 2511.853
$ws.Range("I132").Value = 2511.853
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7535.559
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H134").Value = 16690609
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 16690609
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 50071827
$ws.Range("N134").Value = -50076897
